# Insert a new data row at row 182 ("Hortaliza, Vega Central Mapocho de
# Santiago - Albahaca"): every existing row from 182..240 shifts down by
# one (to 183..241), and the freshly inserted row 182 is populated with a
# new price-report record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 182:240 down to 183:241, carrying formatting along.
$ws.Rows("182:182").Insert()

# Populate the newly-opened row 182 with the new record.
$ws.Cells.Item(182, 1).Value  = 9
$ws.Cells.Item(182, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(182, 3).Value  = "Metropolitana"
$ws.Cells.Item(182, 4).Value  = 44524
$ws.Cells.Item(182, 5).Value  = 13
$ws.Cells.Item(182, 6).Value  = 100112052
$ws.Cells.Item(182, 7).Value  = "Albahaca"
$ws.Cells.Item(182, 8).Value  = "Sin especificar"
$ws.Cells.Item(182, 9).Value  = "Primera"
$ws.Cells.Item(182, 10).Value = 79
$ws.Cells.Item(182, 11).Value = 5000
$ws.Cells.Item(182, 12).Value = 6000
$ws.Cells.Item(182, 13).Value = 5494
$ws.Cells.Item(182, 14).Value = "$/docena de matas"
$ws.Cells.Item(182, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(182, 16).Value = 916
$ws.Cells.Item(182, 17).Value = 6
$ws.Cells.Item(182, 18).Value = "Hortaliza"
